$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 370.26086
$ws.Range("I28").Value = 321.94116
$ws.Range("J28").Value = 507.16666
$ws.Range("K28").Value = 321.94116
$ws.Range("L28").Value = 507.16666
$ws.Range("M28").Value = 163.05884
$ws.Range("N28").Value = -1477.16666
# Row 137
$ws.Range("H137").Value = 3350129.5
$ws.Range("I137").Value = 8548019
$ws.Range("J137").Value = 8628.857
$ws.Range("K137").Value = 25644057
$ws.Range("L137").Value = 25886.571
$ws.Range("M137").Value = -25641507
$ws.Range("N137").Value = -30986.571

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1429.6552
$ws.Range("I2").Value = 1517.76
$ws.Range("J2").Value = 879
$ws.Range("K2").Value = 1517.76
$ws.Range("L2").Value = 879
$ws.Range("M2").Value = -1404.76
$ws.Range("N2").Value = -1105
# Row 45
$ws.Range("H45").Value = 1757.9546
$ws.Range("I45").Value = 1483.7333
$ws.Range("J45").Value = 2345.5715
$ws.Range("K45").Value = 1483.7333
$ws.Range("L45").Value = 2345.5715
$ws.Range("M45").Value = -1106.7333
$ws.Range("N45").Value = -3099.5715
# Row 61
$ws.Range("H61").Value = 2562.5881
$ws.Range("I61").Value = 1837.9166
$ws.Range("J61").Value = 4301.8
$ws.Range("K61").Value = 1837.9166
$ws.Range("L61").Value = 4301.8
$ws.Range("M61").Value = -1625.9166
$ws.Range("N61").Value = -4725.8
# Row 74
$ws.Range("H74").Value = 3077.7273
$ws.Range("I74").Value = 2673
$ws.Range("J74").Value = 4899
$ws.Range("K74").Value = 2673
$ws.Range("L74").Value = 4899
$ws.Range("M74").Value = -1799
$ws.Range("N74").Value = -6647
# Row 77
$ws.Range("H77").Value = 3077.7273
$ws.Range("I77").Value = 2673
$ws.Range("J77").Value = 4899
$ws.Range("K77").Value = 13365
$ws.Range("L77").Value = 24495
$ws.Range("M77").Value = -8997
$ws.Range("N77").Value = -33231
# Row 110
$ws.Range("H110").Value = 1405.5454
$ws.Range("I110").Value = 1328.8334
$ws.Range("J110").Value = 1750.75
$ws.Range("K110").Value = 1328.8334
$ws.Range("L110").Value = 1750.75
$ws.Range("M110").Value = 716.1666
$ws.Range("N110").Value = -5840.75
# Row 116
$ws.Range("H116").Value = 1429.6552
$ws.Range("I116").Value = 1517.76
$ws.Range("J116").Value = 879
$ws.Range("K116").Value = 1517.76
$ws.Range("L116").Value = 879
$ws.Range("M116").Value = 776.24
$ws.Range("N116").Value = -5467
# Row 122
$ws.Range("H122").Value = 1840.125
$ws.Range("I122").Value = 1841.6364
$ws.Range("K122").Value = 5524.9092
$ws.Range("M122").Value = -3074.9092
# Row 136
$ws.Range("H136").Value = 2562.5881
$ws.Range("I136").Value = 1837.9166
$ws.Range("J136").Value = 4301.8
$ws.Range("K136").Value = 5513.7498
$ws.Range("L136").Value = 12905.4
$ws.Range("M136").Value = -2963.7498
$ws.Range("N136").Value = -18005.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1429.6552
$ws.Range("I3").Value = 1517.76
$ws.Range("J3").Value = 879
$ws.Range("K3").Value = 1517.76
$ws.Range("L3").Value = 879
$ws.Range("M3").Value = -1403.76
$ws.Range("N3").Value = -1107
# Row 134
$ws.Range("H134").Value = 3165.982
$ws.Range("I134").Value = 2466.5454
$ws.Range("J134").Value = 3632.2727
$ws.Range("K134").Value = 7399.6362
$ws.Range("L134").Value = 10896.8181
$ws.Range("M134").Value = -4864.6362
$ws.Range("N134").Value = -15966.8181

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2064.4
$ws.Range("I16").Value = 2255.5
$ws.Range("J16").Value = 1300
$ws.Range("K16").Value = 2255.5
$ws.Range("L16").Value = 1300
$ws.Range("M16").Value = -1968.5
$ws.Range("N16").Value = -1874
# Row 31
$ws.Range("H31").Value = 1684.32
$ws.Range("I31").Value = 754.1702
$ws.Range("J31").Value = 2509.17
$ws.Range("K31").Value = 754.1702
$ws.Range("L31").Value = 2509.17
$ws.Range("M31").Value = -459.1702
$ws.Range("N31").Value = -3099.17
# Row 34
$ws.Range("H34").Value = 1684.32
$ws.Range("I34").Value = 754.1702
$ws.Range("J34").Value = 2509.17
$ws.Range("K34").Value = 754.1702
$ws.Range("L34").Value = 2509.17
$ws.Range("M34").Value = -552.1702
$ws.Range("N34").Value = -2913.17
# Row 58
$ws.Range("H58").Value = 1747.8438
$ws.Range("I58").Value = 1430.6666
$ws.Range("J58").Value = 2699.375
$ws.Range("K58").Value = 1430.6666
$ws.Range("L58").Value = 2699.375
$ws.Range("M58").Value = -1227.6666
$ws.Range("N58").Value = -3105.375
# Row 113
$ws.Range("H113").Value = 2064.4
$ws.Range("I113").Value = 2255.5
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 2255.5
$ws.Range("L113").Value = 1300
$ws.Range("M113").Value = -85.5
$ws.Range("N113").Value = -5640
# Row 122
$ws.Range("H122").Value = 240820
$ws.Range("I122").Value = 300725
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 902175
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -899725
$ws.Range("N122").Value = -8500
# Row 132
$ws.Range("H132").Value = 85575.414
$ws.Range("I132").Value = 1407.3
$ws.Range("J132").Value = 205815.58
$ws.Range("K132").Value = 4221.9
$ws.Range("L132").Value = 617446.74
$ws.Range("M132").Value = -1691.9
$ws.Range("N132").Value = -622506.74
# Row 134
$ws.Range("H134").Value = 969968.25
$ws.Range("I134").Value = 587877.8
$ws.Range("J134").Value = 2804002.2
$ws.Range("K134").Value = 1763633.4
$ws.Range("L134").Value = 8412006.600000001
$ws.Range("M134").Value = -1761098.4
$ws.Range("N134").Value = -8417076.600000001
# Row 136
$ws.Range("H136").Value = 1747.8438
$ws.Range("I136").Value = 1430.6666
$ws.Range("J136").Value = 2699.375
$ws.Range("K136").Value = 4291.9998
$ws.Range("L136").Value = 8098.125
$ws.Range("M136").Value = -1741.9998
$ws.Range("N136").Value = -13198.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 2054.6516
$ws.Range("I113").Value = 3413.7942
$ws.Range("J113").Value = 610.5625
$ws.Range("K113").Value = 10241.3826
$ws.Range("L113").Value = 1831.6875
$ws.Range("M113").Value = -8071.382599999999
$ws.Range("N113").Value = -6171.6875
# Row 136
$ws.Range("H136").Value = 29414060
$ws.Range("I136").Value = 55557224
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 166671672
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -166666572
$ws.Range("N136").Value = -19200
# Row 137
$ws.Range("H137").Value = 30307442
$ws.Range("I137").Value = 2757.7856
$ws.Range("J137").Value = 52637210
$ws.Range("K137").Value = 8273.356800000001
$ws.Range("L137").Value = 157911630
$ws.Range("M137").Value = -3173.356800000001
$ws.Range("N137").Value = -157921830

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 58
$ws.Range("H58").Value = 17996.25
$ws.Range("J58").Value = 17996.25
$ws.Range("L58").Value = 17996.25
$ws.Range("N58").Value = -18550.25
# Row 113
$ws.Range("H113").Value = 1977.5555
$ws.Range("J113").Value = 1999.6666
$ws.Range("L113").Value = 1999.6666
$ws.Range("N113").Value = -6339.6666
# Row 122
$ws.Range("H122").Value = 1935.7142
$ws.Range("I122").Value = 1918.1818
$ws.Range("K122").Value = 5754.5454
$ws.Range("M122").Value = -3304.5454
# Row 126
$ws.Range("H126").Value = 10319.077
$ws.Range("I126").Value = 13448.223
$ws.Range("J126").Value = 3278.5
$ws.Range("K126").Value = 40344.669
$ws.Range("L126").Value = 9835.5
$ws.Range("M126").Value = -37874.669
$ws.Range("N126").Value = -14775.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 74
$ws.Range("H74").Value = 45000
$ws.Range("I74").Value = 45000
$ws.Range("K74").Value = 45000
$ws.Range("M74").Value = -44002
# Row 77
$ws.Range("H77").Value = 45000
$ws.Range("I77").Value = 45000
$ws.Range("K77").Value = 135000
$ws.Range("M77").Value = -130008
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 132
$ws.Range("H132").Value = 2233.1538
$ws.Range("I132").Value = 1710.6981
$ws.Range("J132").Value = 4540.6665
$ws.Range("K132").Value = 5132.094300000001
$ws.Range("L132").Value = 13621.9995
$ws.Range("M132").Value = -2602.094300000001
$ws.Range("N132").Value = -18681.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 468.13333
$ws.Range("I113").Value = 413.22223
$ws.Range("K113").Value = 1239.66669
$ws.Range("M113").Value = 930.33331
# Row 122
$ws.Range("H122").Value = 2858403
$ws.Range("I122").Value = 2858403
$ws.Range("K122").Value = 8575209
$ws.Range("M122").Value = -8572759
# Row 136
$ws.Range("H136").Value = 1667878.1
$ws.Range("I136").Value = 2593466.2
$ws.Range("J136").Value = 1819.4
$ws.Range("K136").Value = 7780398.600000001
$ws.Range("L136").Value = 5458.200000000001
$ws.Range("M136").Value = -7777848.600000001
$ws.Range("N136").Value = -10558.2
